$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows describing the JSON form read by NodoTracker.
$ws.Range("A13").Value = "JSON"

$ws.Range("A14").Value = "CANT NODOS"
$ws.Range("C14").Value = "??"

$ws.Range("A15").Value = "IP SIGUIENTE"
$ws.Range("A16").Value = "PORT SIGUIENTE"

$ws.Range("C15").Value = "-"
$ws.Range("C16").Value = 27015

# Underline the "HASH 1" label to call it out, same as an author highlighting it.
$ws.Range("B5").Font.Underline = $true

# Select B5 and I14 (matches the saved selection state in the file),
# leaving I14 -- the new form area -- as the active cell.
$excel.Union($ws.Range("B5"), $ws.Range("I14")).Select()
$ws.Range("I14").Select()
